# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.388.44"
$ws.Range("E2").Value = "  +2.95%  "
# Row 3
$ws.Range("D3").Value = "'2.639.65"
$ws.Range("E3").Value = "  +1.36%  "
# Row 5
$ws.Range("D5").Value = "'604.66"
$ws.Range("E5").Value = "  +2.85%  "
# Row 6
$ws.Range("D6").Value = "'155.72"
$ws.Range("E6").Value = "  +4.28%  "
# Row 7
$ws.Range("D7").Value = "'0.999"
# Row 8
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +0.77%  "
# Row 9
$ws.Range("E9").Value = "  +8.85%  "
# Row 10
$ws.Range("D10").Value = "'0.409"
$ws.Range("E10").Value = "  +5.96%  "
# Row 11
$ws.Range("E11").Value = "  +0.10%  "
# Row 12
$ws.Range("E12").Value = "  +3.08%  "
# Row 13
$ws.Range("D13").Value = "'29.36"
$ws.Range("E13").Value = "  +6.58%  "
# Row 14
$ws.Range("E14").Value = "  +22.85%  "
# Row 15
$ws.Range("D15").Value = "'3.109.25"
$ws.Range("E15").Value = "  +1.28%  "
# Row 16
$ws.Range("D16").Value = "'65.194.24"
$ws.Range("E16").Value = "  +2.97%  "
# Row 17
$ws.Range("D17").Value = "'2.649.92"
$ws.Range("E17").Value = "  +1.95%  "
# Row 18
$ws.Range("D18").Value = "'12.61"
$ws.Range("E18").Value = "  +4.39%  "
# Row 19
$ws.Range("E19").Value = "  +5.52%  "
# Row 20
$ws.Range("D20").Value = "'359.98"
# Row 21
$ws.Range("D21").Value = "'7.37"
$ws.Range("E21").Value = "  +7.72%  "
# Row 22
$ws.Range("E22").Value = "  +0.08%  "
# Row 23
$ws.Range("D23").Value = "'68.89"
$ws.Range("E23").Value = "  +3.68%  "
# Row 24
$ws.Range("E24").Value = "  -3.48%  "
# Row 25
$ws.Range("E25").Value = "  +2.49%  "
# Row 26
$ws.Range("D26").Value = "'1.66"
$ws.Range("E26").Value = "  +0.30%  "
# Row 27
$ws.Range("D27").Value = "'8.25"
$ws.Range("E27").Value = "  +0.81%  "
# Row 28
$ws.Range("E28").Value = "  +2.95%  "
# Row 29
$ws.Range("D29").Value = "'0.0₃0963"
$ws.Range("E29").Value = "  +14.11%  "
# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.19"
$ws.Range("E30").Value = "  +8.42%  "
# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.10%  "
# Row 32
$ws.Range("D32").Value = "'527.75"
$ws.Range("E32").Value = "  -6.49%  "
# Row 33
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  +2.78%  "
# Row 34
$ws.Range("D34").Value = "'5.52"
$ws.Range("E34").Value = "  +4.51%  "
# Row 35
$ws.Range("D35").Value = "'6.39"
$ws.Range("E35").Value = "  +5.14%  "
# Row 36
$ws.Range("E36").Value = "  +4.04%  "
# Row 37
$ws.Range("D37").Value = "'20.55"
$ws.Range("E37").Value = "  +5.85%  "
# Row 38
$ws.Range("E38").Value = "  +5.47%  "
# Row 39
$ws.Range("D39").Value = "'162.37"
$ws.Range("E39").Value = "  -2.02%  "
# Row 40
$ws.Range("E40").Value = "  -0.03%  "
# Row 41
$ws.Range("E41").Value = "  +0.07%  "
# Row 42
$ws.Range("D42").Value = "'42.85"
$ws.Range("E42").Value = "  +8.33%  "
# Row 43
$ws.Range("D43").Value = "'165.46"
$ws.Range("E43").Value = "  -0.13%  "
# Row 44
$ws.Range("D44").Value = "'4.18"
$ws.Range("E44").Value = "  +5.00%  "
# Row 45
$ws.Range("D45").Value = "'0.0621"
$ws.Range("E45").Value = "  +7.28%  "
# Row 46
$ws.Range("D46").Value = "'23.46"
$ws.Range("E46").Value = "  +3.03%  "
# Row 47
$ws.Range("D47").Value = "'2.24"
$ws.Range("E47").Value = "  +6.47%  "
# Row 48
$ws.Range("E48").Value = "  +7.89%  "
# Row 49
$ws.Range("D49").Value = "'0.656"
$ws.Range("E49").Value = "  +4.01%  "
# Row 50
$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = "  +2.68%  "
# Row 51
$ws.Range("D51").Value = "'19.72"
$ws.Range("E51").Value = "  +3.29%  "
